# The deck ships two theme parts:
#   theme1.xml ("Integral")    -> used by the (only) slide master
#   theme2.xml ("Office Theme")-> used by the notes master
# The authored edit swaps the two themes' content, so the slide master
# picks up the plain "Office Theme" palette (and the notes master would
# pick up "Integral" -- not reachable from this object model, see below).
#
# The PowerPoint OM doesn't expose the raw theme XML, but
# ThemeColorScheme.Colors(i).RGB is a live round-trip onto the theme's
# <a:clrScheme> entries (order = MsoThemeColorSchemeIndex: dk1, lt1, dk2,
# lt2, accent1-6, hlink, folHlink). Drive every slide's ThemeColorScheme
# to the "Office Theme" default palette to realize that half of the swap.

$p = $ppt.ActivePresentation

# Office Theme default palette, in MsoThemeColorSchemeIndex order (1..12),
# expressed as OLE (BGR) RGB() integers: dk1, lt1, dk2, lt2,
# accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $officeThemeColors.Length; $i++) {
    $tcs.Colors($i).RGB = $officeThemeColors[$i - 1]
}
